$d = $word.ActiveDocument

# =====================================================================
# Section 1: "Data de entrega: 30/03" paragraph ->
#            "Datas de entrega: 31/08 e 01/09"
# =====================================================================

$p4 = $d.Paragraphs(4)
$start4 = $p4.Range.Start

# Insert "s" right after "Data" to get "Datas"
$insS = $d.Range($start4 + 4, $start4 + 4)
$insS.Text = "s"

# Replace the date "30/03" with "31/08 e 01/09"
$p4b = $d.Paragraphs(4)
$findDate = $d.Range($p4b.Range.Start, $p4b.Range.End)
$null = $findDate.Find.Execute("30/03", $false, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
$findDate.Text = "31/08 e 01/09"

# =====================================================================
# Section 2: OBS paragraph - update student/group counts and move the
#            "_GoBack" bookmark to sit right before "Os grupos podem..."
# =====================================================================

$p21 = $d.Paragraphs(21)
$pStart = $p21.Range.Start
$pEnd = $p21.Range.End

# Replace the whole "pessoas. Como são 27 ... diferentes. " run of text
# in one shot (keeps things simple, formatting is fixed up afterwards).
$seg = $d.Range($pStart, $pEnd)
$null = $seg.Find.Execute( `
    "pessoas. Como são 27 alunos, serão 9 grupos de 3 pessoas. Os grupos podem ser compostos por integrantes de turmas diferentes. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg.Text = "pessoas. Como são 30 alunos no total, serão 10 grupos de 3 pessoas. Os grupos podem ser compostos por integrantes de turmas diferentes. "

# "pessoas." (without the trailing space) should not be bold.
$p21b = $d.Paragraphs(21)
$findPessoas = $d.Range($p21b.Range.Start, $p21b.Range.End)
$null = $findPessoas.Find.Execute("pessoas. Como", $false, $false, $false, `
                                   $false, $false, $true, 1, $false, "", 0)
$pessoasStart = $findPessoas.Start
$pessoasRange = $d.Range($pessoasStart, $pessoasStart + 8)
$pessoasRange.Font.Bold = 0

# The closing " Os grupos podem ser compostos por integrantes de turmas
# diferentes. " tail should not be bold either.
$p21c = $d.Paragraphs(21)
$findTail = $d.Range($p21c.Range.Start, $p21c.Range.End)
$null = $findTail.Find.Execute( `
    "pessoas. Os grupos podem ser compostos por integrantes de turmas diferentes. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $findTail.Start + 8
$tailRange = $d.Range($splitPoint, $findTail.End)
$tailRange.Font.Bold = 0

# Move the "_GoBack" bookmark here (this automatically removes it from
# its old spot at the end of the "Data(s) de entrega" paragraph).
$bmRange = $d.Range($splitPoint, $splitPoint)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
